$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.290.96"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.227.98"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$cell = $ws.Range("D5")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "244.69"
$cell.Style = $s
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  +1.13%  "
$cell = $ws.Range("D7")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "74.22"
$cell.Style = $s
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  +0.11%  "
$cell = $ws.Range("D9")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.616"
$cell.Style = $s
$ws.Range("E9").Value = "  -0.19%  "
$cell = $ws.Range("D10")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "43.11"
$cell.Style = $s
$ws.Range("E10").Value = "  +4.58%  "
$cell = $ws.Range("D11")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0959"
$cell.Style = $s
$ws.Range("E11").Value = "  +1.80%  "
$cell = $ws.Range("D12")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$cell.Style = $s
$ws.Range("E12").Value = "  +1.20%  "
$cell = $ws.Range("D13")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.102"
$cell.Style = $s
$ws.Range("E13").Value = "  -0.45%  "
$cell = $ws.Range("D14")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.43"
$cell.Style = $s
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "2.236.42"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "42.140.17"
$ws.Range("E17").Value = "  +0.59%  "
$cell = $ws.Range("D18")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000110"
$cell.Style = $s
$ws.Range("E18").Value = "  +12.25%  "
$ws.Range("E19").Value = "  +1.98%  "
$cell = $ws.Range("D20")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "72.05"
$cell.Style = $s
$ws.Range("E20").Value = "  +0.77%  "
$cell = $ws.Range("D21")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.17"
$cell.Style = $s
$ws.Range("E21").Value = "  +38.64%  "
$cell = $ws.Range("D22")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "231.01"
$cell.Style = $s
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  -4.43%  "
$cell = $ws.Range("D24")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.78"
$cell.Style = $s
$ws.Range("E24").Value = "  +5.62%  "
$ws.Range("E25").Value = "  +0.13%  "
$cell = $ws.Range("D26")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.66"
$cell.Style = $s
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  +1.24%  "
$cell = $ws.Range("D28")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.Style = $s
$ws.Range("E28").Value = "  +3.17%  "
$cell = $ws.Range("D29")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "166.98"
$cell.Style = $s
$ws.Range("E29").Value = "  -0.89%  "
$cell = $ws.Range("D30")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.02"
$cell.Style = $s
$ws.Range("E30").Value = "  +2.59%  "
$cell = $ws.Range("D31")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.82"
$cell.Style = $s
$ws.Range("E31").Value = "  +18.46%  "
$cell = $ws.Range("D32")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0805"
$cell.Style = $s
$ws.Range("E32").Value = "  -2.15%  "
$cell = $ws.Range("D33")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.118"
$cell.Style = $s
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  +0.02%  "
$cell = $ws.Range("D35")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.54"
$cell.Style = $s
$ws.Range("E35").Value = "  -7.93%  "
$cell = $ws.Range("D36")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.40"
$cell.Style = $s
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("E37").Value = "  +2.69%  "
$cell = $ws.Range("D38")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.24"
$cell.Style = $s
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("E39").Value = "  +0.04%  "
$cell = $ws.Range("D40")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = $s
$ws.Range("E40").Value = "  -3.18%  "
$cell = $ws.Range("D41")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.29"
$cell.Style = $s
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  +2.23%  "
$cell = $ws.Range("D44")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "104.62"
$cell.Style = $s
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("E45").Value = "  +3.10%  "
$cell = $ws.Range("D46")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.993"
$cell.Style = $s
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +6.39%  "
$ws.Range("E48").Value = "  +0.33%  "
$cell = $ws.Range("D49")
$s = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.18"
$cell.Style = $s
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("E51").Value = "  -0.87%  "
